$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 ("Secondary") education-level percentages were recomputed
# (age-standardization removed on the age-group breakdown) -- update
# columns B:L with the new values, keeping them as text like the rest
# of the table.
$newValues = @{
    "B10" = "0.29"
    "C10" = "0.49"
    "D10" = "0.37"
    "E10" = "0.51"
    "F10" = "0.42"
    "G10" = "0.58"
    "H10" = "0.5"
    "I10" = "0.63"
    "J10" = "0.65"
    "K10" = "0.66"
    "L10" = "0.61"
}

foreach ($addr in $newValues.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $newValues[$addr]
}
